$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E1 value
$ws.Range("E1").Value = 0.5649590492248535

# Add new row 2 with data
$ws.Range("A2").Value = "ResultsA3.csv"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 42
$ws.Range("E2").Value = 29.99599480628967
$ws.Range("F2").Value = 123
$ws.Range("G2").Value = 124462
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 125
